$d = $word.ActiveDocument

# --- Change 1: split the "Known Issues" Heading1 run, wrapping "Known" in
#     proofErr gramStart/gramEnd markers (mirrors Word's live grammar-check
#     run-splitting behaviour). ---
$headingTarget = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Known Issues*") {
        $headingTarget = $p
        break
    }
}

if ($headingTarget -ne $null) {
    $headingXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:t>Known</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> Issues</w:t></w:r></w:p>'
    $headingTarget.Range.InsertXML($headingXml) | Out-Null
}

# --- Change 2: add a new sub-bullet after "Port from Electron to Tauri
#     (eventually)" at the same list level (ilvl=1, numId=3). ---
$portTarget = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Port from Electron to Tauri*") {
        $portTarget = $p
        break
    }
}

if ($portTarget -ne $null) {
    $portTarget.Range.InsertParagraphAfter()

    $paras = @($d.Paragraphs)
    for ($i = 0; $i -lt $paras.Count; $i++) {
        if ($paras[$i].Range.Text -like "Port from Electron to Tauri*") {
            $newPara = $paras[$i + 1]
            break
        }
    }

    $newPara.Range.Text = "Dealing with images better, every song loads its cover into memory which leads to duplicates. (idk)"
}

Write-Output "done"
